$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($rangeAddr, $val)
    $rng = $ws.Range($rangeAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '29.126.08'
Set-TextValue "E2" '  -0.33%  '

Set-TextValue "D3" '1.853.98'
Set-TextValue "E3" '  +0.04%  '

Set-TextValue "D4" '0.9997'
Set-TextValue "E4" '  +0.00%  '

$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue "D5" '0.6902'
Set-TextValue "E5" '  -1.08%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue "D6" '237.44'
Set-TextValue "E6" '  -0.52%  '

Set-TextValue "E7" '  +0.04%  '

Set-TextValue "D8" '0.07740'
Set-TextValue "E8" '  +1.72%  '

Set-TextValue "D9" '0.3033'
Set-TextValue "E9" '  -1.26%  '

Set-TextValue "D10" '22.95'
Set-TextValue "E10" '  -2.86%  '

Set-TextValue "D11" '0.08074'
Set-TextValue "E11" '  -0.16%  '

Set-TextValue "D12" '1.842.02'
Set-TextValue "E12" '  -0.54%  '

Set-TextValue "D13" '0.7189'
Set-TextValue "E13" '  -0.93%  '

Set-TextValue "D14" '5.155'
Set-TextValue "E14" '  -0.63%  '

Set-TextValue "D15" '89.34'
Set-TextValue "E15" '  +0.26%  '

Set-TextValue "D16" '29.134.44'
Set-TextValue "E16" '  -0.67%  '

Set-TextValue "D17" '5.716'
Set-TextValue "E17" '  -2.66%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D18" '0.000007751'
Set-TextValue "E18" '  +0.35%  '

$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D19" '13.20'
Set-TextValue "E19" '  +0.62%  '

Set-TextValue "D20" '234.69'
Set-TextValue "E20" '  -2.92%  '

Set-TextValue "D21" '0.9994'
Set-TextValue "E21" '  -0.05%  '

Set-TextValue "D22" '2.110.09'
Set-TextValue "E22" '  -0.14%  '

Set-TextValue "E23" '  +0.04%  '

Set-TextValue "D24" '7.451'
Set-TextValue "E24" '  -2.27%  '

Set-TextValue "D25" '8.985'
Set-TextValue "E25" '  -0.77%  '

Set-TextValue "D26" '161.18'
Set-TextValue "E26" '  -0.40%  '

Set-TextValue "D27" '0.1429'
Set-TextValue "E27" '  -2.32%  '

Set-TextValue "D28" '18.00'
Set-TextValue "E28" '  -0.34%  '

Set-TextValue "D29" '1.939'
Set-TextValue "E29" '  +0.18%  '

Set-TextValue "D30" '1.407'
Set-TextValue "E30" '  +0.44%  '

Set-TextValue "D31" '4.477'
Set-TextValue "E31" '  +0.81%  '

Set-TextValue "D32" '1.484'
Set-TextValue "E32" '  -1.23%  '

Set-TextValue "D33" '3.999'
Set-TextValue "E33" '  -1.16%  '

Set-TextValue "D34" '0.05175'
Set-TextValue "E34" '  -1.56%  '

Set-TextValue "D35" '1.169'
Set-TextValue "E35" '  -2.00%  '

Set-TextValue "D36" '0.7057'
Set-TextValue "E36" '  -0.83%  '

Set-TextValue "D37" '1.001'

Set-TextValue "D38" '2.656'
Set-TextValue "E38" '  -0.20%  '

Set-TextValue "D39" '0.01847'
Set-TextValue "E39" '  -0.84%  '

Set-TextValue "D40" '2.716'
Set-TextValue "E40" '  +1.42%  '

Set-TextValue "D41" '0.9316'
Set-TextValue "E41" '  +0.69%  '

Set-TextValue "D42" '1.106.31'
Set-TextValue "E42" '  +5.74%  '

Set-TextValue "D43" '0.4266'
Set-TextValue "E43" '  -0.80%  '

Set-TextValue "D44" '5.874'
Set-TextValue "E44" '  -1.26%  '

Set-TextValue "D45" '70.33'
Set-TextValue "E45" '  +0.94%  '

Set-TextValue "E46" '  +0.03%  '

Set-TextValue "D47" '102.75'
Set-TextValue "E47" '  +0.39%  '

Set-TextValue "D48" '1.785'
Set-TextValue "E48" '  +2.53%  '

Set-TextValue "D49" '2.006.72'
Set-TextValue "E49" '  -0.20%  '

Set-TextValue "D50" '9.108'
Set-TextValue "E50" '  -1.66%  '

Set-TextValue "D51" '6.962'
Set-TextValue "E51" '  -3.88%  '
